$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "MIE(5.7637555488771, 4.247257195142881, -6.659484915168087, 9.258208568900521)"
$ws.Range("C2").Value = "NIG(0.7659046416502266, 0.501112600305588, 5.936722188046554, 4.979771815534203)"
$ws.Range("D2").Value = "NCT(2.9650191419577445, 1.6158338413730648, -0.9670632119792966, 2.437354425456652)"
$ws.Range("E2").Value = "NIG(1.049084254482982, 0.807297829251432, 4.85164001064399, 5.026121155240573)"
